$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells for the new "suavizamiento exponencial" (alpha / ypredicho) block ---
$ws.Range("K1").Value = "alpha"
$ws.Range("L1").Value = 0.99
$ws.Range("I1").Value = "ypredicho"

# --- h=1 labels for rows 2..31 (one-step-ahead smoothing) ---
$ws.Range("H2:H31").Value = "h=1"

# --- ypredicho (simple exponential smoothing) ---
# I2 is the seed forecast = first observed value
$ws.Range("I2").Formula = "=B2"
# I3 = alpha*B2 + (1-alpha)*I2  (first recursive step, not part of the shared group)
$ws.Range("I3").Formula = '=$L$1*B2+(1-$L$1)*I2'
# I4:I31 share the same relative pattern: alpha*B(row-1) + (1-alpha)*I(row-1)
$ws.Range("I4:I31").Formula = '=$L$1*B3+(1-$L$1)*I3'

# --- Holt-style flat forecasts for future horizons h=2..4 ---
$ws.Range("H32").Value = "h=2"
$ws.Range("I32").Formula = '=$L$1*I31+(1-$L$1)*I31'

$ws.Range("H33").Value = "h=3"
$ws.Range("I33").Formula = '=$L$1*I32+(1-$L$1)*I32'

$ws.Range("H34").Value = "h=4"
$ws.Range("I34").Formula = '=$L$1*I33+(1-$L$1)*I33'

# --- View / navigation state ---
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("K10").Select()
